# Remove form_id from remaining forms
#
# The "settings" sheet has columns: form_title, form_id, version, style, namespaces
# This removes the form_id column (and its associated comment / shared string),
# shifting version/style/namespaces left by one column.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- settings sheet: shift header comments left (B1 <- C1 <- D1 <- E1), drop E1 ---
# (Comment for B1/form_id is replaced by the comment that belongs to the next
# column over, and so on, so that after the form_id column is deleted the
# comments still describe the correct header.)
$versionComment   = $settings.Range("C1").Comment.Text()
$styleComment     = $settings.Range("D1").Comment.Text()
$namespacesComment = $settings.Range("E1").Comment.Text()

$settings.Range("B1").Comment.Text($versionComment)
$settings.Range("C1").Comment.Text($styleComment)
$settings.Range("D1").Comment.Text($namespacesComment)
$settings.Range("E1").Comment.Delete()

# --- settings sheet: delete the form_id column (B) entirely ---
$settings.Range("B1:B2").EntireColumn.Delete()

# --- settings sheet: restore the cursor/selection to B1 ---
$settings.Activate()
$settings.Range("B1").Select()

# --- survey sheet: update the remembered selection for the bottom-right pane ---
$survey.Activate()
$survey.Range("A11").Select()

# --- survey sheet: consolidate the conditional-formatting range for column C ---
# (C26:C9998, C2:C24 and C25 cover exactly the same cells as C2:C9998; the
# rule/formula/priority/dxf stay the same, only the stored range is tidied up.)
$fcs = $survey.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Formula1 -eq '=AND(ISBLANK(C2),NOT(OR(ISBLANK($A2),$A2="calculate")))') {
        $fc.ModifyAppliesToRange($survey.Range("C2:C9998"))
        break
    }
}

# keep "survey" as the active sheet/tab, same as in the original workbook
$survey.Activate()
$survey.Range("A11").Select()
